$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.163.93"
$ws.Range("E2").Value = "'  +0.35%  "
$ws.Range("D3").Value = "'1.657.30"
$ws.Range("E3").Value = "'  -0.22%  "
$ws.Range("E4").Value = "'  -0.22%  "
$ws.Range("D5").Value = "'215.42"
$ws.Range("E5").Value = "'  +3.48%  "
$ws.Range("D6").Value = "'0.5236"
$ws.Range("E6").Value = "'  +1.21%  "
$ws.Range("E7").Value = "'  -0.17%  "
$ws.Range("D8").Value = "'0.2628"
$ws.Range("E8").Value = "'  +1.76%  "
$ws.Range("D9").Value = "'0.06393"
$ws.Range("E9").Value = "'  +1.43%  "
$ws.Range("D10").Value = "'20.88"
$ws.Range("E10").Value = "'  -0.31%  "
$ws.Range("D11").Value = "'0.07767"
$ws.Range("E11").Value = "'  +3.07%  "
$ws.Range("D12").Value = "'1.653.81"
$ws.Range("E12").Value = "'  -0.57%  "
$ws.Range("D13").Value = "'4.461"
$ws.Range("E13").Value = "'  +1.39%  "
$ws.Range("D14").Value = "'1.882.06"
$ws.Range("E14").Value = "'  -0.42%  "
$ws.Range("D15").Value = "'0.5526"
$ws.Range("E15").Value = "'  +2.71%  "
$ws.Range("D16").Value = "'0.0₅8275"
$ws.Range("E16").Value = "'  +4.11%  "
$ws.Range("D17").Value = "'65.17"
$ws.Range("E17").Value = "'  -1.49%  "
$ws.Range("D18").Value = "'26.179.32"
$ws.Range("E18").Value = "'  +0.36%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "'  -0.15%  "
$ws.Range("D20").Value = "'4.759"
$ws.Range("E20").Value = "'  +1.19%  "
$ws.Range("D21").Value = "'190.43"
$ws.Range("E21").Value = "'  +1.76%  "
$ws.Range("D22").Value = "'10.31"
$ws.Range("D23").Value = "'6.374"
$ws.Range("E23").Value = "'  +2.89%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "'  -0.23%  "
$ws.Range("D25").Value = "'143.18"
$ws.Range("E25").Value = "'  -3.47%  "
$ws.Range("D26").Value = "'0.1252"
$ws.Range("E26").Value = "'  +3.13%  "
$ws.Range("D27").Value = "'7.412"
$ws.Range("E27").Value = "'  +0.44%  "
$ws.Range("D28").Value = "'15.99"
$ws.Range("E28").Value = "'  +2.21%  "
$ws.Range("D29").Value = "'1.431"
$ws.Range("E29").Value = "'  +2.88%  "
$ws.Range("D30").Value = "'0.06109"
$ws.Range("E30").Value = "'  +1.88%  "
$ws.Range("D31").Value = "'1.266"
$ws.Range("E31").Value = "'  +0.28%  "
$ws.Range("D32").Value = "'3.515"
$ws.Range("E32").Value = "'  +1.37%  "
$ws.Range("D33").Value = "'3.422"
$ws.Range("E33").Value = "'  +0.71%  "
$ws.Range("D34").Value = "'1.663"
$ws.Range("E34").Value = "'  +1.55%  "
$ws.Range("D35").Value = "'1.0000"
$ws.Range("E35").Value = "'  +1.49%  "
$ws.Range("D36").Value = "'2.400"
$ws.Range("E36").Value = "'  +0.56%  "
$ws.Range("D37").Value = "'2.759"
$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("D38").Value = "'0.5662"
$ws.Range("E38").Value = "'  -3.74%  "
$ws.Range("E39").Value = "'  +0.67%  "
$ws.Range("D40").Value = "'5.908"
$ws.Range("E40").Value = "'  -1.04%  "
$ws.Range("D41").Value = "'0.8550"
$ws.Range("E41").Value = "'  +0.94%  "
$ws.Range("E42").Value = "'  -0.16%  "
$ws.Range("D43").Value = "'1.031.95"
$ws.Range("E43").Value = "'  -6.58%  "
$ws.Range("D44").Value = "'99.62"
$ws.Range("E44").Value = "'  -0.27%  "
$ws.Range("D45").Value = "'1.805.02"
$ws.Range("E45").Value = "'  -0.68%  "
$ws.Range("D46").Value = "'0.0₈107"
$ws.Range("E46").Value = "'  -1.50%  "
$ws.Range("D47").Value = "'56.14"
$ws.Range("E47").Value = "'  +1.91%  "
$ws.Range("E48").Value = "'  +0.28%  "
$ws.Range("D49").Value = "'8.111"
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'5.977"
$ws.Range("E51").Value = "'  +2.09%  "
